$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.046.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.636.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.32"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5244"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2597"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06293"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.71"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07661"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.632.47"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.416"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.859.66"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5531"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8220"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.95"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.031.31"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.695"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "187.85"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.18"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.154"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.22"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1216"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.410"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.394"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05963"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.25%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.436"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.413"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.644"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9851"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.04%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5670"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8502"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.758"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.69%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.035.34"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -6.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.24"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.785.24"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈108"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.66"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.040"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4214"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.71%  "
